$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date (column C) for rows 2-13 from 2023-11-03 (45233)
# to 2023-11-13 (45243), keeping the existing date format/style intact.
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 3).Value2 = 45243
}
